# Update the "Förändrad" (Changed) date column (C) for rows 2 through 28
# from 45185 (2023-09-16) to 45204 (2023-10-05).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C28").Value = 45204

Write-Host "Updated C2:C28 to 45204"
